$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the duplicate AIST row (old row 9, "Michihiro Ohta").
#    This shifts the old row 10 (NIMS) up to become row 9.
$ws.Rows("9").Delete()

# 2. Update text content that changed between the old and new versions.

# Row 3 (KIMS): author name updated / expanded
$ws.Range("E3").Value = "Eun-Ae CHOI; Seungzeon HAN"

# Row 4 (KAIST): author name capitalization fixed, collaboration note added
$ws.Range("E4").Value = "Seunghwa RYU"
$ws.Range("F4").Value = "Advanced TE Device"

# Row 5 (DLR): authors reordered, collaboration text tweaked
$ws.Range("E5").Value = "Johannes de Boor; Pawel Ziolkowski; Eckhard Mueller"
$ws.Range("F5").Value = "TGM characterization, silicide TE mater."

# Row 8 (AIST): merge in the second author that used to be its own row, add collaboration text
$ws.Range("E8").Value = "Michihiro Ohta; Atsushi YAMAMOTO"
$ws.Range("F8").Value = "TGM characterization"

# Row 9 (was row 10, NIMS): trim author name, add collaboration text
$ws.Range("E9").Value = "Yukari Katsura"
$ws.Range("F9").Value = "TE data and efficiency map"

# 3. Apply word-wrap formatting to the "who" and "Collaboration" columns for all data rows.
$ws.Range("E2:F9").WrapText = $true

# 4. Adjust row heights to match the new wrapped content.
$ws.Rows("3").RowHeight = 33.75
$ws.Rows("5").RowHeight = 33.75
$ws.Rows("8").RowHeight = 33.75
$ws.Rows("9").RowHeight = 33.75

# 5. Reset the view: clear the frozen/scrolled top-left cell and select F10.
$ws.Range("A1").Select() | Out-Null
$ws.Range("F10").Select() | Out-Null
